$d = $word.ActiveDocument

# Update the date line (first paragraph)
$d.Content.Find.Execute("2025-01-17 Friday", $true, $false, $false, $false, $false, $true, 1, $false, "2025-01-18 Saturday", 2) | Out-Null

# Update the table of math expressions, cell by cell (row-major, 20 rows x 5 cols)
$t = $d.Tables.Item(1)
$newValues = @(
    "77-22=",
    "2+33=",
    "83-32=",
    "80-9=",
    "90-15=",
    "5+29=",
    "33-30=",
    "18+53=",
    "97-8=",
    "0+42=",
    "92-21=",
    "88-81=",
    "72-16=",
    "1+34=",
    "83-36=",
    "3+39=",
    "22+2=",
    "71-32=",
    "20+48=",
    "78-65=",
    "6+38=",
    "66-49=",
    "72+23=",
    "32+41=",
    "81-35=",
    "11+72=",
    "23-12=",
    "25+51=",
    "87-5=",
    "76-26=",
    "47-31=",
    "74-57=",
    "27-10=",
    "62-40=",
    "45+33=",
    "78-7=",
    "25+39=",
    "2-1=",
    "9+82=",
    "89-76=",
    "59+38=",
    "7+1=",
    "50-15=",
    "8+42=",
    "52-26=",
    "24+3=",
    "52-5=",
    "52-7=",
    "93-0=",
    "0+15=",
    "50-24=",
    "58-0=",
    "48+34=",
    "90-29=",
    "2+7=",
    "57+14=",
    "74-23=",
    "63+3=",
    "39+0=",
    "22+73=",
    "83-78=",
    "1+91=",
    "32+11=",
    "8+11=",
    "3+13=",
    "78+16=",
    "26-4=",
    "70+24=",
    "99-5=",
    "27+64=",
    "49+0=",
    "64-30=",
    "21-13=",
    "45+12=",
    "19+25=",
    "72+22=",
    "1+65=",
    "96-58=",
    "35+57=",
    "79-1=",
    "6+76=",
    "91-13=",
    "23-3=",
    "82-30=",
    "79-52=",
    "53-27=",
    "43-21=",
    "64-19=",
    "16+41=",
    "77-71=",
    "59-37=",
    "35+45=",
    "66-15=",
    "2+30=",
    "55+27=",
    "31+38=",
    "77-20=",
    "52+3=",
    "23+41=",
    "43+26="
)

$idx = 0
for ($r = 1; $r -le $t.Rows.Count; $r++) {
    for ($c = 1; $c -le $t.Columns.Count; $c++) {
        $cell = $t.Cell($r, $c)
        $rng = $cell.Range
        $rng.End = $rng.End - 1
        $rng.Text = $newValues[$idx]
        $idx = $idx + 1
    }
}
